$d = $word.ActiveDocument

# The document currently starts with:
#   Para 1 (Heading1): "Questionnaire - Draft"
#   Para 2 (FirstParagraph): "Notes" (italic) + ":"
#
# We need to insert a new FirstParagraph paragraph with the text
# "Advertised as: Study on Online News and Political Attitudes" right
# before the "Notes:" paragraph, and change the "Notes:" paragraph's
# style from FirstParagraph to BodyText (keeping "Notes" italic).

$notesPara = $d.Paragraphs(2)
$rng = $notesPara.Range.Duplicate
$rng.Collapse(1)

$introText = "Advertised as: Study on Online News and Political Attitudes"
$insertStart = $rng.Start

# Insert the new text right before "Notes" - InsertBefore creates it as its
# own run (no inherited direct character formatting), unlike assigning
# Range.Text on a freshly split empty paragraph.
$rng.InsertBefore($introText)

# Split the paragraph between our new text and "Notes" so the new text
# becomes its own paragraph (keeping the original FirstParagraph style),
# and the "Notes:" text remains in the paragraph that follows.
$splitPos = $insertStart + $introText.Length
$splitRng = $d.Range($splitPos, $splitPos)
$splitRng.InsertParagraphBefore()

# The original "Notes:" paragraph is now the 3rd paragraph; restyle it.
$notesPara2 = $d.Paragraphs(3)
$notesRunRange = $d.Range($notesPara2.Range.Start, $notesPara2.Range.Start + 5)
$notesPara2.Style = "BodyText"
# Re-apply the italic direct formatting on "Notes" that the style change
# clears, so it stays italic like in the original.
$notesRunRange.Font.Italic = -1

Write-Host "done"
